$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before row 125, shifting existing rows 125-171 down to 128-174.
$ws.Rows.Item(125).Resize(3).Insert()

# Data for the three newly inserted rows (125, 126, 127).
$newRows = @(
    @{ H = "Sin especificar"; I = "Banquete"; J = 500; K = 1500; L = 1600; M = 1546; N = "`$/kilo"; O = "Provincia de Linares"; P = 1546; Q = 1 },
    @{ H = "Sin especificar"; I = "Primera";  J = 540; K = 1100; L = 1200; M = 1156; N = "`$/kilo"; O = "Provincia de Linares"; P = 1156; Q = 1 },
    @{ H = "Sin especificar"; I = "Segunda";  J = 550; K = 900;  L = 1000; M = 958;  N = "`$/kilo"; O = "Provincia de Linares"; P = 958;  Q = 1 }
)

$rowIdx = 125
foreach ($row in $newRows) {
    $ws.Cells.Item($rowIdx, 1).Value = 6
    $ws.Cells.Item($rowIdx, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
    $ws.Cells.Item($rowIdx, 3).Value = "Metropolitana"
    $ws.Cells.Item($rowIdx, 4).Value = (Get-Date -Year 2023 -Month 10 -Day 3 -Hour 0 -Minute 0 -Second 0)
    $ws.Cells.Item($rowIdx, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $ws.Cells.Item($rowIdx, 5).Value = 13
    $ws.Cells.Item($rowIdx, 6).Value = 300000000
    $ws.Cells.Item($rowIdx, 7).Value = "Espárragos"
    $ws.Cells.Item($rowIdx, 8).Value = $row.H
    $ws.Cells.Item($rowIdx, 9).Value = $row.I
    $ws.Cells.Item($rowIdx, 10).Value = $row.J
    $ws.Cells.Item($rowIdx, 11).Value = $row.K
    $ws.Cells.Item($rowIdx, 12).Value = $row.L
    $ws.Cells.Item($rowIdx, 13).Value = $row.M
    $ws.Cells.Item($rowIdx, 14).Value = $row.N
    $ws.Cells.Item($rowIdx, 15).Value = $row.O
    $ws.Cells.Item($rowIdx, 16).Value = $row.P
    $ws.Cells.Item($rowIdx, 17).Value = $row.Q
    $ws.Cells.Item($rowIdx, 18).Value = "Hortaliza"
    $rowIdx++
}
